$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.779013514518738
$ws.Range("B1").Value = 2.283915042877197
$ws.Range("C1").Value = 2.463985204696655
$ws.Range("D1").Value = 6.540905952453613
$ws.Range("E1").Value = 0.7693328261375427
